$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
  # Row 13: The Hexster Runoff | Maple Picatrix
  $ws.Range("H13").Value = 8000.3076
  $ws.Range("J13").Value = 14751
  $ws.Range("L13").Value = 14751
  $ws.Range("N13").Value = -15089
  # Row 106: Making Your Mark | Enchanted Palladium Ink
  $ws.Range("H106").Value = 3153.9473
  $ws.Range("I106").Value = 2561.8667
  $ws.Range("K106").Value = 2561.8667
  $ws.Range("M106").Value = -1930.8667
  # Row 132: Fast-forwarding Flora | Growth Formula Lambda
  $ws.Range("H132").Value = 3579.3635
  $ws.Range("I132").Value = 3828.8948
  $ws.Range("K132").Value = 11486.6844
  $ws.Range("M132").Value = -8956.6844

$ws = $wb.Worksheets.Item("ARM")
  # Row 32: Ingot We Trust | Steel Ingot
  $ws.Range("H32").Value = 4183.05
  $ws.Range("I32").Value = 2617.5686
  $ws.Range("K32").Value = 2617.5686
  $ws.Range("M32").Value = -2330.5686
  # Row 61: Dealing with the Tough Stuff | Cobalt Ingot
  $ws.Range("H61").Value = 24394072
  $ws.Range("I61").Value = 24394072
  $ws.Range("K61").Value = 24394072
  $ws.Range("M61").Value = -24393860
  # Row 102: Smells of Rich Tama-hagane | Tama-hagane Ingot
  $ws.Range("H102").Value = 693.7778
  $ws.Range("I102").Value = 705.17645
  $ws.Range("J102").Value = 500
  $ws.Range("K102").Value = 705.17645
  $ws.Range("L102").Value = 500
  $ws.Range("M102").Value = 916.82355
  $ws.Range("N102").Value = -3744
  # Row 122: Haste for High Durium | High Durium Nugget
  $ws.Range("H122").Value = 2247.3845
  $ws.Range("I122").Value = 2338.0908
  $ws.Range("J122").Value = 1748.5
  $ws.Range("K122").Value = 7014.2724
  $ws.Range("L122").Value = 5245.5
  $ws.Range("M122").Value = -4564.2724
  $ws.Range("N122").Value = -10145.5
  # Row 127: Once and for Alchemy | Bismuth Alembic
  $ws.Range("H127").Value = 100000
  $ws.Range("J127").Value = 100000
  $ws.Range("L127").Value = 100000
  $ws.Range("N127").Value = -109920
  # Row 136: Metal with Mettle | Cobalt Tungsten Ingot
  $ws.Range("H136").Value = 24394072
  $ws.Range("I136").Value = 24394072
  $ws.Range("K136").Value = 73182216
  $ws.Range("M136").Value = -73179666

$ws = $wb.Worksheets.Item("BSM")
  # Row 76: Keep Up with the Mechanics | Titanium-barreled Arquebus
  $ws.Range("H76").Value = 19066.666
  $ws.Range("J76").Value = 19066.666
  $ws.Range("L76").Value = 19066.666
  $ws.Range("N76").Value = -19696.666
  # Row 79: Unconventional Weaponry (L) | Titanium-barreled Arquebus
  $ws.Range("H79").Value = 19066.666
  $ws.Range("J79").Value = 19066.666
  $ws.Range("L79").Value = 19066.666
  $ws.Range("N79").Value = -21250.666
  # Row 105: Ingot to Wing It | Molybdenum Ingot
  $ws.Range("H105").Value = 2852.25
  $ws.Range("I105").Value = 2148
  $ws.Range("K105").Value = 2148
  $ws.Range("M105").Value = -401
  # Row 107: The Gold Experience | Deepgold Nugget
  $ws.Range("H107").Value = 94744.55
  $ws.Range("I107").Value = 1663.1666
  $ws.Range("J107").Value = 206442.2
  $ws.Range("K107").Value = 1663.1666
  $ws.Range("L107").Value = 206442.2
  $ws.Range("M107").Value = 256.8334
  $ws.Range("N107").Value = -210282.2
  # Row 134: Ruthenium Supremium | Ruthenium Ingot
  $ws.Range("H134").Value = 10002219
  $ws.Range("I134").Value = 10871216
  $ws.Range("K134").Value = 32613648
  $ws.Range("M134").Value = -32611113

$ws = $wb.Worksheets.Item("CRP")
  # Row 19: Shielding Sales | Square Ash Shield
  $ws.Range("H19").Value = 2345.125
  $ws.Range("I19").Value = 3043
  $ws.Range("J19").Value = 1182
  $ws.Range("K19").Value = 3043
  $ws.Range("L19").Value = 1182
  $ws.Range("M19").Value = -2873
  $ws.Range("N19").Value = -1522
  # Row 24: What You Need | Square Ash Shield
  $ws.Range("H24").Value = 2345.125
  $ws.Range("I24").Value = 3043
  $ws.Range("J24").Value = 1182
  $ws.Range("K24").Value = 3043
  $ws.Range("L24").Value = 1182
  $ws.Range("M24").Value = -2873
  $ws.Range("N24").Value = -1522
  # Row 31: Wall Not Found | Walnut Lumber
  $ws.Range("H31").Value = 8234.657999999999
  $ws.Range("I31").Value = 1951.6666
  $ws.Range("K31").Value = 1951.6666
  $ws.Range("M31").Value = -1656.6666
  # Row 34: Armoires of the Rich and Famous | Walnut Lumber
  $ws.Range("H34").Value = 8234.657999999999
  $ws.Range("I34").Value = 1951.6666
  $ws.Range("K34").Value = 1951.6666
  $ws.Range("M34").Value = -1749.6666
  # Row 74: License to Heal | Dark Chestnut Rod
  $ws.Range("H74").Value = 0
  $ws.Range("J74").Value = 0
  $ws.Range("L74").Value = 0
  $ws.Range("N74").Value = $null
  # Row 77: Purified Polyrhythm (L) | Dark Chestnut Rod
  $ws.Range("H77").Value = 0
  $ws.Range("J77").Value = 0
  $ws.Range("L77").Value = 0
  $ws.Range("N77").Value = $null
  # Row 132: Hull Lotta Damage | Ginseng Lumber
  $ws.Range("H132").Value = 125000970
  $ws.Range("I132").Value = 200000960
  $ws.Range("K132").Value = 600002880
  $ws.Range("M132").Value = -600000350
  # Row 134: Wood You Be Quiet | Ceiba Lumber
  $ws.Range("H134").Value = 25000886
  $ws.Range("I134").Value = 25000886
  $ws.Range("K134").Value = 75002658
  $ws.Range("M134").Value = -75000123
  # Row 141: No Greater Treasure | Claro Walnut Necklace of Gathering
  $ws.Range("H141").Value = 230653.16
  $ws.Range("J141").Value = 267135.53
  $ws.Range("L141").Value = 267135.53
  $ws.Range("N141").Value = -277495.53

$ws = $wb.Worksheets.Item("CUL")
  # Row 22: A Total Nut Job | Walnut Bread
  $ws.Range("H22").Value = 0
  $ws.Range("I22").Value = 0
  $ws.Range("K22").Value = 0
  $ws.Range("M22").Value = $null
  # Row 27: Brain Food | Walnut Bread
  $ws.Range("H27").Value = 0
  $ws.Range("I27").Value = 0
  $ws.Range("K27").Value = 0
  $ws.Range("M27").Value = $null
  # Row 60: Drinking to Your Health | Mulled Tea
  $ws.Range("H60").Value = 5802.0586
  $ws.Range("I60").Value = 283.75
  $ws.Range("K60").Value = 851.25
  $ws.Range("M60").Value = -600.25
  # Row 82: Persuasion of a Higher Power | Baked Pipira Pira
  $ws.Range("H82").Value = 0
  $ws.Range("I82").Value = 0
  $ws.Range("K82").Value = 0
  $ws.Range("M82").Value = $null
  # Row 85: Loaves and Fishes (L) | Baked Pipira Pira
  $ws.Range("H85").Value = 0
  $ws.Range("I85").Value = 0
  $ws.Range("K85").Value = 0
  $ws.Range("M85").Value = $null
  # Row 92: Oh No Udon | Gyr Abanian Flour
  $ws.Range("H92").Value = 611
  $ws.Range("I92").Value = 611
  $ws.Range("K92").Value = 1833
  $ws.Range("M92").Value = -585
  # Row 98: Sweet Kiss of Death | Rice Vinegar
  $ws.Range("H98").Value = 766.5
  $ws.Range("I98").Value = 1177.6666
  $ws.Range("J98").Value = 355.33334
  $ws.Range("K98").Value = 3532.9998
  $ws.Range("L98").Value = 1066.00002
  $ws.Range("M98").Value = -2034.9998
  $ws.Range("N98").Value = -4062.00002
  # Row 123: Topping Up the Pot | Zurek
  $ws.Range("H123").Value = 6966.6665
  # Row 129: Comfort Food | Yakow Moussaka
  $ws.Range("H129").Value = 2480.182
  $ws.Range("I129").Value = 430.5
  $ws.Range("J129").Value = 4939.8
  $ws.Range("K129").Value = 1291.5
  $ws.Range("L129").Value = 14819.4
  $ws.Range("M129").Value = 3708.5
  $ws.Range("N129").Value = -24819.4
  # Row 131: The Mountain Steeped | Tsai tou Vounou
  $ws.Range("H131").Value = 2176.923
  $ws.Range("I131").Value = 1862
  $ws.Range("K131").Value = 5586
  $ws.Range("M131").Value = -546
  # Row 137: Creative Chocolate | Gateau au Chocolat
  $ws.Range("H137").Value = 6251921
  $ws.Range("I137").Value = 9092549
  $ws.Range("K137").Value = 27277647
  $ws.Range("M137").Value = -27272547

$ws = $wb.Worksheets.Item("GSM")
  # Row 113: Copious Crystal Cannons | Manasilver Nugget
  $ws.Range("H113").Value = 66222
  $ws.Range("I113").Value = 93336.63
  $ws.Range("J113").Value = 6569.8
  $ws.Range("K113").Value = 93336.63
  $ws.Range("L113").Value = 6569.8
  $ws.Range("M113").Value = -91166.63
  $ws.Range("N113").Value = -10909.8
  # Row 119: Bulking Up | Dwarven Mythril Rapier
  $ws.Range("H119").Value = 53760.668
  $ws.Range("J119").Value = 53760.668
  $ws.Range("L119").Value = 53760.668
  $ws.Range("N119").Value = -63436.668
  # Row 132: On Board for Lar | Lar Ingot
  $ws.Range("H132").Value = 7354431
  $ws.Range("I132").Value = 10418056
  $ws.Range("J132").Value = 1731.4
  $ws.Range("K132").Value = 31254168
  $ws.Range("L132").Value = 5194.200000000001
  $ws.Range("M132").Value = -31251638
  $ws.Range("N132").Value = -10254.2

$ws = $wb.Worksheets.Item("LTW")
  # Row 40: Best Served Toad | Toad Leather
  $ws.Range("H40").Value = 4203.5835
  $ws.Range("I40").Value = 4203.5835
  $ws.Range("K40").Value = 4203.5835
  $ws.Range("M40").Value = -4067.5835
  # Row 46: Supply Side Logic | Boar Leather
  $ws.Range("H46").Value = 739.2857
  $ws.Range("J46").Value = 499
  $ws.Range("L46").Value = 499
  $ws.Range("N46").Value = -875
  # Row 82: Trainin' the Neck | Dragon Leather
  $ws.Range("H82").Value = 1118.5714
  $ws.Range("I82").Value = 966.3
  $ws.Range("K82").Value = 966.3
  $ws.Range("M82").Value = -605.3
  # Row 85: Training Is Only Skintight (L) | Dragon Leather
  $ws.Range("H85").Value = 1118.5714
  $ws.Range("I85").Value = 966.3
  $ws.Range("K85").Value = 966.3
  $ws.Range("M85").Value = 281.7
  # Row 131: For What Was Gleaned | Ophiotauroskin Wristband of Gathering
  $ws.Range("H131").Value = 100000
  $ws.Range("J131").Value = 100000
  $ws.Range("L131").Value = 100000
  $ws.Range("N131").Value = -110080
  # Row 132: Tenets of Tanning | Silver Lobo Leather
  $ws.Range("H132").Value = 36933524
  $ws.Range("I132").Value = 40011156
  $ws.Range("K132").Value = 120033468
  $ws.Range("M132").Value = -120030938

$ws = $wb.Worksheets.Item("WVR")
  # Row 33: I'll Be Your Wailer Today | Velveteen Wedge Cap of Gathering
  $ws.Range("H33").Value = 5840
  $ws.Range("I33").Value = 4749.5
  $ws.Range("J33").Value = 8021
  $ws.Range("K33").Value = 4749.5
  $ws.Range("L33").Value = 8021
  $ws.Range("M33").Value = -4499.5
  $ws.Range("N33").Value = -8521
  # Row 36: Put a Lid on It | Velveteen Wedge Cap of Gathering
  $ws.Range("H36").Value = 5840
  $ws.Range("I36").Value = 4749.5
  $ws.Range("J36").Value = 8021
  $ws.Range("K36").Value = 4749.5
  $ws.Range("L36").Value = 8021
  $ws.Range("M36").Value = -4499.5
  $ws.Range("N36").Value = -8521
  # Row 122: Heavy Armoire | Dark Hempen Cloth
  $ws.Range("H122").Value = 5231.8823
  $ws.Range("I122").Value = 4836.1333
  $ws.Range("K122").Value = 14508.3999
  $ws.Range("M122").Value = -12058.3999
  # Row 132: Comfy Cabins | Snow Cotton Cloth
  $ws.Range("H132").Value = 11114141
  $ws.Range("I132").Value = 12198214
  $ws.Range("J132").Value = 2399.25
  $ws.Range("K132").Value = 36594642
  $ws.Range("L132").Value = 7197.75
  $ws.Range("M132").Value = -36592112
  $ws.Range("N132").Value = -12257.75
  # Row 136: Weaving the Envelope | Sarcenet Cloth
  $ws.Range("H136").Value = 31252002
  $ws.Range("I136").Value = 41667452
  $ws.Range("K136").Value = 125002356
  $ws.Range("M136").Value = -124999806
